$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2018.9445
$ws.Range("J17").Value = 1973.0588
$ws.Range("L17").Value = 5919.1764
$ws.Range("N17").Value = -6255.1764
$ws.Range("H51").Value = 7139.222
$ws.Range("J51").Value = 6191.25
$ws.Range("L51").Value = 6191.25
$ws.Range("N51").Value = -7159.25
$ws.Range("H64").Value = 1000000000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 1000000000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H106").Value = 2742.318
$ws.Range("I106").Value = 2701.55
$ws.Range("K106").Value = 2701.55
$ws.Range("M106").Value = -2070.55
$ws.Range("H111").Value = 775.6667
$ws.Range("I111").Value = 647.5
$ws.Range("K111").Value = 1942.5
$ws.Range("M111").Value = 1124.5
$ws.Range("H116").Value = 3999.5
$ws.Range("J116").Value = 3999.5
$ws.Range("L116").Value = 3999.5
$ws.Range("N116").Value = -10883.5
$ws.Range("H138").Value = 309428.12
$ws.Range("I138").Value = 6149.846
$ws.Range("J138").Value = 366567.53
$ws.Range("K138").Value = 18449.538
$ws.Range("L138").Value = 1099702.59
$ws.Range("M138").Value = -13309.538
$ws.Range("N138").Value = -1109982.59

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4096.7437
$ws.Range("I32").Value = 3629
$ws.Range("J32").Value = 12750
$ws.Range("K32").Value = 3629
$ws.Range("L32").Value = 12750
$ws.Range("M32").Value = -3342
$ws.Range("N32").Value = -13324
$ws.Range("H45").Value = 19935.75
$ws.Range("I45").Value = 28095.688
$ws.Range("K45").Value = 28095.688
$ws.Range("M45").Value = -27718.688
$ws.Range("H61").Value = 3357.0852
$ws.Range("I61").Value = 1479.2778
$ws.Range("J61").Value = 9502.637000000001
$ws.Range("K61").Value = 1479.2778
$ws.Range("L61").Value = 9502.637000000001
$ws.Range("M61").Value = -1267.2778
$ws.Range("N61").Value = -9926.637000000001
$ws.Range("H122").Value = 2514.4285
$ws.Range("I122").Value = 2372.465
$ws.Range("K122").Value = 7117.395
$ws.Range("M122").Value = -4667.395
$ws.Range("H130").Value = 81995.336
$ws.Range("J130").Value = 81995.336
$ws.Range("L130").Value = 81995.336
$ws.Range("N130").Value = -92035.336
$ws.Range("H132").Value = 2069.4465
$ws.Range("I132").Value = 1727.7174
$ws.Range("K132").Value = 5183.1522
$ws.Range("M132").Value = -2653.1522
$ws.Range("H136").Value = 3357.0852
$ws.Range("I136").Value = 1479.2778
$ws.Range("J136").Value = 9502.637000000001
$ws.Range("K136").Value = 4437.8334
$ws.Range("L136").Value = 28507.911
$ws.Range("M136").Value = -1887.8334
$ws.Range("N136").Value = -33607.911

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 34727160
$ws.Range("I20").Value = 41671816
$ws.Range("J20").Value = 3871
$ws.Range("K20").Value = 41671816
$ws.Range("L20").Value = 3871
$ws.Range("M20").Value = -41671569
$ws.Range("N20").Value = -4365
$ws.Range("H26").Value = 2973
$ws.Range("I26").Value = 2973
$ws.Range("K26").Value = 2973
$ws.Range("M26").Value = -2681
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H86").Value = 3478.509
$ws.Range("I86").Value = 3176.0715
$ws.Range("K86").Value = 3176.0715
$ws.Range("M86").Value = -2053.0715
$ws.Range("H89").Value = 3478.509
$ws.Range("I89").Value = 3176.0715
$ws.Range("K89").Value = 15880.3575
$ws.Range("M89").Value = -10264.3575
$ws.Range("H107").Value = 1673.875
$ws.Range("I107").Value = 1443.7059
$ws.Range("J107").Value = 2232.8572
$ws.Range("K107").Value = 1443.7059
$ws.Range("L107").Value = 2232.8572
$ws.Range("M107").Value = 476.2941000000001
$ws.Range("N107").Value = -6072.8572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2195.5
$ws.Range("I105").Value = 1567.7273
$ws.Range("K105").Value = 1567.7273
$ws.Range("M105").Value = 179.2727
$ws.Range("H122").Value = 2174.2273
$ws.Range("I122").Value = 2262.2
$ws.Range("K122").Value = 6786.599999999999
$ws.Range("M122").Value = -4336.599999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1171.95
$ws.Range("I2").Value = 94.63636
$ws.Range("K2").Value = 567.81816
$ws.Range("M2").Value = -454.81816
$ws.Range("H33").Value = 527473.3
$ws.Range("I33").Value = 833448.5
$ws.Range("J33").Value = 2944.4285
$ws.Range("K33").Value = 5000691
$ws.Range("L33").Value = 17666.571
$ws.Range("M33").Value = -5000408
$ws.Range("N33").Value = -18232.571
$ws.Range("H63").Value = 3237.8572
$ws.Range("I63").Value = 1666.75
$ws.Range("J63").Value = 5332.6665
$ws.Range("K63").Value = 5000.25
$ws.Range("L63").Value = 15997.9995
$ws.Range("M63").Value = -4251.25
$ws.Range("N63").Value = -17495.9995
$ws.Range("H64").Value = 2963.5454
$ws.Range("I64").Value = 300
$ws.Range("J64").Value = 3229.9
$ws.Range("K64").Value = 900
$ws.Range("L64").Value = 9689.700000000001
$ws.Range("M64").Value = -630
$ws.Range("N64").Value = -10229.7
$ws.Range("H66").Value = 3237.8572
$ws.Range("I66").Value = 1666.75
$ws.Range("J66").Value = 5332.6665
$ws.Range("K66").Value = 15000.75
$ws.Range("L66").Value = 47993.9985
$ws.Range("M66").Value = -11256.75
$ws.Range("N66").Value = -55481.9985
$ws.Range("H67").Value = 2963.5454
$ws.Range("I67").Value = 300
$ws.Range("J67").Value = 3229.9
$ws.Range("K67").Value = 900
$ws.Range("L67").Value = 9689.700000000001
$ws.Range("M67").Value = 36
$ws.Range("N67").Value = -11561.7
$ws.Range("H107").Value = 830.0526
$ws.Range("I107").Value = 595.4
$ws.Range("K107").Value = 1786.2
$ws.Range("M107").Value = 133.8000000000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14348084
$ws.Range("I70").Value = 19312576
$ws.Range("J70").Value = 6218.1113
$ws.Range("K70").Value = 19312576
$ws.Range("L70").Value = 6218.1113
$ws.Range("M70").Value = -19312306
$ws.Range("N70").Value = -6758.1113
$ws.Range("H73").Value = 14348084
$ws.Range("I73").Value = 19312576
$ws.Range("J73").Value = 6218.1113
$ws.Range("K73").Value = 19312576
$ws.Range("L73").Value = 6218.1113
$ws.Range("M73").Value = -19311640
$ws.Range("N73").Value = -8090.1113
$ws.Range("H122").Value = 2483446.2
$ws.Range("I122").Value = 2566127.5
$ws.Range("K122").Value = 7698382.5
$ws.Range("M122").Value = -7695932.5
$ws.Range("H126").Value = 9895.5
$ws.Range("I126").Value = 3467.2
$ws.Range("J126").Value = 20609.334
$ws.Range("K126").Value = 10401.6
$ws.Range("L126").Value = 61828.00199999999
$ws.Range("M126").Value = -7931.599999999999
$ws.Range("N126").Value = -66768.00199999999
$ws.Range("H132").Value = 3609.093
$ws.Range("I132").Value = 2764.3572
$ws.Range("J132").Value = 5185.933
$ws.Range("K132").Value = 8293.071599999999
$ws.Range("L132").Value = 15557.799
$ws.Range("M132").Value = -5763.071599999999
$ws.Range("N132").Value = -20617.799

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2157.682
$ws.Range("I61").Value = 2157.7058
$ws.Range("J61").Value = 2157.6
$ws.Range("K61").Value = 2157.7058
$ws.Range("L61").Value = 2157.6
$ws.Range("M61").Value = -1955.7058
$ws.Range("N61").Value = -2561.6
$ws.Range("H93").Value = 1651
$ws.Range("I93").Value = 1651
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1651
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -403
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 4680.2
$ws.Range("I100").Value = 4384.8945
$ws.Range("J100").Value = 5615.3335
$ws.Range("K100").Value = 4384.8945
$ws.Range("L100").Value = 5615.3335
$ws.Range("M100").Value = -3843.8945
$ws.Range("N100").Value = -6697.3335
$ws.Range("H113").Value = 2157.682
$ws.Range("I113").Value = 2157.7058
$ws.Range("J113").Value = 2157.6
$ws.Range("K113").Value = 2157.7058
$ws.Range("L113").Value = 2157.6
$ws.Range("M113").Value = 12.29419999999982
$ws.Range("N113").Value = -6497.6
$ws.Range("H132").Value = 4744.643
$ws.Range("I132").Value = 2903.8462
$ws.Range("K132").Value = 8711.5386
$ws.Range("M132").Value = -6181.5386
$ws.Range("H136").Value = 8340.091
$ws.Range("I136").Value = 8999.333000000001
$ws.Range("J136").Value = 8092.875
$ws.Range("K136").Value = 26997.999
$ws.Range("L136").Value = 24278.625
$ws.Range("M136").Value = -24447.999
$ws.Range("N136").Value = -29378.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64127.57
$ws.Range("J46").Value = 64127.57
$ws.Range("L46").Value = 64127.57
$ws.Range("N46").Value = -64589.57
$ws.Range("H81").Value = 4465.6665
$ws.Range("I81").Value = 3359.7693
$ws.Range("J81").Value = 5772.636
$ws.Range("K81").Value = 6719.5386
$ws.Range("L81").Value = 11545.272
$ws.Range("M81").Value = -5658.5386
$ws.Range("N81").Value = -13667.272
$ws.Range("H84").Value = 4465.6665
$ws.Range("I84").Value = 3359.7693
$ws.Range("J84").Value = 5772.636
$ws.Range("K84").Value = 33597.693
$ws.Range("L84").Value = 57726.36
$ws.Range("M84").Value = -28293.693
$ws.Range("N84").Value = -68334.36
$ws.Range("H94").Value = 8888
$ws.Range("J94").Value = 8888
$ws.Range("L94").Value = 8888
$ws.Range("N94").Value = -10690
$ws.Range("H134").Value = 64127.57
$ws.Range("J134").Value = 64127.57
$ws.Range("L134").Value = 192382.71
$ws.Range("N134").Value = -197452.71
$ws.Range("H140").Value = 83987.57000000001
$ws.Range("J140").Value = 83987.57000000001
$ws.Range("L140").Value = 83987.57000000001
$ws.Range("N140").Value = -94347.57000000001
